$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before row 703 to shift existing data down
$ws.Rows("703:711").Insert()

# Ensure the date/id columns keep plain text (not auto-converted to dates/numbers)
$ws.Range("B703:B711").NumberFormat = "@"
$ws.Range("C703:C711").NumberFormat = "@"

# Populate the newly inserted rows with historical data (2019-11-18 .. 2019-11-28)
$ws.Range("A703").Value = 1574035200
$ws.Range("B703").Value = "2019-11-18"
$ws.Range("C703").Value = "5277"
$ws.Range("D703").Value = "FPGROUP"
$ws.Range("E703").Value = 0.545
$ws.Range("F703").Value = 0.555
$ws.Range("G703").Value = 0.53
$ws.Range("H703").Value = 0.545
$ws.Range("I703").Value = 15209300

$ws.Range("A704").Value = 1574121600
$ws.Range("B704").Value = "2019-11-19"
$ws.Range("C704").Value = "5277"
$ws.Range("D704").Value = "FPGROUP"
$ws.Range("E704").Value = 0.545
$ws.Range("F704").Value = 0.55
$ws.Range("G704").Value = 0.535
$ws.Range("H704").Value = 0.545
$ws.Range("I704").Value = 12472900

$ws.Range("A705").Value = 1574208000
$ws.Range("B705").Value = "2019-11-20"
$ws.Range("C705").Value = "5277"
$ws.Range("D705").Value = "FPGROUP"
$ws.Range("E705").Value = 0.55
$ws.Range("F705").Value = 0.58
$ws.Range("G705").Value = 0.545
$ws.Range("H705").Value = 0.58
$ws.Range("I705").Value = 30667800

$ws.Range("A706").Value = 1574294400
$ws.Range("B706").Value = "2019-11-21"
$ws.Range("C706").Value = "5277"
$ws.Range("D706").Value = "FPGROUP"
$ws.Range("E706").Value = 0.58
$ws.Range("F706").Value = 0.585
$ws.Range("G706").Value = 0.57
$ws.Range("H706").Value = 0.58
$ws.Range("I706").Value = 11679100

$ws.Range("A707").Value = 1574380800
$ws.Range("B707").Value = "2019-11-22"
$ws.Range("C707").Value = "5277"
$ws.Range("D707").Value = "FPGROUP"
$ws.Range("E707").Value = 0.58
$ws.Range("F707").Value = 0.58
$ws.Range("G707").Value = 0.57
$ws.Range("H707").Value = 0.575
$ws.Range("I707").Value = 4244600

$ws.Range("A708").Value = 1574640000
$ws.Range("B708").Value = "2019-11-25"
$ws.Range("C708").Value = "5277"
$ws.Range("D708").Value = "FPGROUP"
$ws.Range("E708").Value = 0.575
$ws.Range("F708").Value = 0.575
$ws.Range("G708").Value = 0.5600000000000001
$ws.Range("H708").Value = 0.5649999999999999
$ws.Range("I708").Value = 6340300

$ws.Range("A709").Value = 1574726400
$ws.Range("B709").Value = "2019-11-26"
$ws.Range("C709").Value = "5277"
$ws.Range("D709").Value = "FPGROUP"
$ws.Range("E709").Value = 0.5649999999999999
$ws.Range("F709").Value = 0.57
$ws.Range("G709").Value = 0.555
$ws.Range("H709").Value = 0.555
$ws.Range("I709").Value = 6356200

$ws.Range("A710").Value = 1574812800
$ws.Range("B710").Value = "2019-11-27"
$ws.Range("C710").Value = "5277"
$ws.Range("D710").Value = "FPGROUP"
$ws.Range("E710").Value = 0.555
$ws.Range("F710").Value = 0.58
$ws.Range("G710").Value = 0.555
$ws.Range("H710").Value = 0.575
$ws.Range("I710").Value = 11499700

$ws.Range("A711").Value = 1574899200
$ws.Range("B711").Value = "2019-11-28"
$ws.Range("C711").Value = "5277"
$ws.Range("D711").Value = "FPGROUP"
$ws.Range("E711").Value = 0.575
$ws.Range("F711").Value = 0.615
$ws.Range("G711").Value = 0.57
$ws.Range("H711").Value = 0.6
$ws.Range("I711").Value = 28204400
